# Update NATMI LR-pair sheet with refreshed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Ligand-expressing cells" count (2, was 1) is identical across all data rows (2-5),
# which cascades into detection rate / average / total expression of the ligand column,
# and the ligand is shared by rows 2-5.
$ligandExpressingCells   = 2
$ligandDetectionRate     = 0.6666666666666666
$ligandAvgExpr           = 0.8077336666666667
$ligandTotalExpr         = 2.423201

$rows = @(2, 3, 4, 5)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value  = $ligandExpressingCells   # E: Ligand-expressing cells
    $ws.Cells.Item($r, 6).Value  = $ligandDetectionRate     # F: Ligand detection rate
    $ws.Cells.Item($r, 7).Value  = $ligandAvgExpr           # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $ligandTotalExpr         # H: Ligand total expression value
}

# Receptor average/total expression values refreshed with new TPM data (column M/N),
# and their derived specificities (O/P) plus edge weights/specificities (Q/R/S/T)
# recomputed accordingly for each target cluster row.
$receptorData = @{
    2 = @{ M = 0.74396;            N = 2.23188;
           O = 0.006259003216804254; P = 0.006259003216804255;
           Q = 0.6009215386533333;   R = 5.40829384788;
           S = 0.006259003216804254; T = 0.006259003216804255 }
    3 = @{ M = 88.14978533333333;  N = 264.449356;
           O = 0.7416121699579786;  P = 0.7416121699579786;
           Q = 71.20154932317288;   R = 640.813943908556;
           S = 0.7416121699579786;  T = 0.7416121699579786 }
    4 = @{ M = 29.76859933333333;  N = 89.305798;
           O = 0.2504459365921425;  P = 0.2504459365921425;
           Q = 24.04509989104422;   R = 216.405899019398;
           S = 0.2504459365921425;  T = 0.2504459365921425 }
    5 = @{ M = 0.2000323333333334; N = 0.6000970000000001;
           O = 0.00168289023307462; P = 0.00168289023307462;
           Q = 0.1615728500552223;  R = 1.454155650497;
           S = 0.00168289023307462; T = 0.00168289023307462 }
}

foreach ($r in $rows) {
    $vals = $receptorData[$r]
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Receptor derived specificity of total expression value
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $vals.R   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $vals.S   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $vals.T   # T: Edge total expression derived specificity
}
